$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Address column (E) with "x" suffixed test values
$ws.Range("E2").Value = "D109x"
$ws.Range("E3").Value = "D110x"
$ws.Range("E4").Value = "D111x"
$ws.Range("E5").Value = "D112x"

# Update Verified * column (N) to "Yes" for rows 2 and 3
$ws.Range("N2").Value = "Yes"
$ws.Range("N3").Value = "Yes"

# Update Update Only * column (O) to "Yes" for all data rows
$ws.Range("O2").Value = "Yes"
$ws.Range("O3").Value = "Yes"
$ws.Range("O4").Value = "Yes"
$ws.Range("O5").Value = "Yes"

# Update selected cell to O5 to match the saved selection state
$ws.Range("O5").Select()
